# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
#
# Refreshes the Belgium MSME summary figures with more precise decimal
# values (density, employment %, enterprises %, and value-added % for
# Micro / SMEs / MSMEs). The source values were stored as text, so the
# replacements are written as text as well (leading apostrophe forces
# text entry, same as typing it directly into Excel), and the style is
# reset back to "Normal" afterwards so the text-quote-prefix formatting
# introduced by the apostrophe doesn't linger on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
Set-TextValue "B13" "49.98"
Set-TextValue "C13" "2.87"
Set-TextValue "D13" "52.85"

# Employment (% of total): Micro / SMEs / MSMEs
Set-TextValue "B14" "32.68"
Set-TextValue "C14" "36.36"
Set-TextValue "D14" "69.04"

# Enterprises (% of total): Micro / SMEs / MSMEs
Set-TextValue "B16" "94.43"
Set-TextValue "C16" "5.42"
Set-TextValue "D16" "99.85"

# Value added to the economy (% of total): Micro / SMEs / MSMEs
Set-TextValue "B20" "22.97"
Set-TextValue "C20" "38.87"
Set-TextValue "D20" "61.84"
